# Adding J=0 and rerunning analysis
# Insert a new row at row 11 (J=0), pushing the existing J=-0.1 .. J=-0.9
# rows down from 11-19 to 12-20, then populate the new row with the
# freshly computed results for J=0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 11:19 down to 12:20, leaving row 11 empty for the new data.
$ws.Rows.Item(11).Insert()

# Populate the new J=0 row.
$ws.Range("A11").Value = 35
$ws.Range("B11").Value = 0.1
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0.483
$ws.Range("E11").Value = 0.022
$ws.Range("F11").Value = -0.062
$ws.Range("G11").Value = 0.128
$ws.Range("H11").Value = 0.476
$ws.Range("I11").Value = 0.026

# Match the author's final selection in the workbook.
$ws.Range("J11").Select()
